# Updated symbol list on Wed Dec 28 20:44:47 UTC 2022 with GitHub Actions
#
# Applies the per-cell text updates to the "cryptos" sheet:
#  - Refreshed Price values for several existing rows
#  - Rows 18-24 and 41-43 shift: coin identity (Coin/Link/Volume) moves down
#    one slot, with rows 24 (One) and 41-43 (KickToken/BKEXToken/CEJI)
#    ending up re-ordered/re-labelled, each carrying its own new Price
#  - Row 48 (BOLO) gets a refreshed Price and loses the "Bestin24h" suffix
#
# All affected cells are stored as literal text (t="inlineStr" in the
# original file) - NOT numbers - even though most of the Price column looks
# numeric. Directly assigning a numeric-looking string via Range.Value makes
# Excel coerce it to a real number, which would change the cell's stored
# type. To keep these as text (matching the source workbook), values are
# staged in a scratch cell formatted as Text and copied over with
# PasteSpecial (values only), which preserves the Text type. The scratch
# row is removed again once all values have been copied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row just below the sheet's used range, used as scratch space so we don't
# disturb any real data; removed (EntireRow.Delete) once we're done.
$scratchRow = 52
$scratch = $ws.Cells.Item($scratchRow, 1)
$scratch.NumberFormat = "@"

function Set-TextValue([int]$row, [int]$col, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
}

# ----- Simple Price-only refreshes -----
Set-TextValue 2  4 "244.01"
Set-TextValue 3  4 "23.79"
Set-TextValue 4  4 "5.245"
Set-TextValue 6  4 "6.464"
Set-TextValue 7  4 "3.230"
Set-TextValue 9  4 "0.8786"
Set-TextValue 10 4 "0.1394"
Set-TextValue 11 4 "0.07088"
Set-TextValue 12 4 "0.03187"
Set-TextValue 14 4 "0.09330"
Set-TextValue 15 4 "3.844"
Set-TextValue 16 4 "0.001542"
Set-TextValue 17 4 "0.04701"

# ----- Row 18: One -> TigerCash -----
$ws.Cells.Item(18, 2).Value = "TigerCash"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue 18 4 "0.006191"
$ws.Cells.Item(18, 5).Value = "17TigerCashTCH"

# ----- Row 19: TigerCash -> BitKan -----
$ws.Cells.Item(19, 2).Value = "BitKan"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue 19 4 "0.001258"
$ws.Cells.Item(19, 5).Value = "18BitKanKAN"

# ----- Row 20: BitKan -> HotbitToken -----
$ws.Cells.Item(20, 2).Value = "HotbitToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue 20 4 "0.004070"
$ws.Cells.Item(20, 5).Value = "19HotbitTokenHTB"

# ----- Row 21: HotbitToken -> NitroEx -----
$ws.Cells.Item(21, 2).Value = "NitroEx"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue 21 4 "0.00008703"
$ws.Cells.Item(21, 5).Value = "20NitroExNTX"

# ----- Row 22: NitroEx -> LEO -----
$ws.Cells.Item(22, 2).Value = "LEO"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue 22 4 "3.545"
$ws.Cells.Item(22, 5).Value = "21LEOLEO"

# ----- Row 23: LEO -> BTSEToken -----
$ws.Cells.Item(23, 2).Value = "BTSEToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue 23 4 "2.154"
$ws.Cells.Item(23, 5).Value = "22BTSETokenBTSE"

# ----- Row 24: BTSEToken -> One -----
$ws.Cells.Item(24, 2).Value = "One"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 24 4 "0.01032"
$ws.Cells.Item(24, 5).Value = "23OneONEBestin24h"

# ----- Simple Price-only refreshes (continued) -----
Set-TextValue 40 4 "0.03787"

# ----- Row 41: BKEXToken -> KickToken -----
$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 41 4 "0.006242"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"

# ----- Row 42: CEJI -> BKEXToken -----
$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1050"
$ws.Cells.Item(42, 5).Value = "41BKEXTokenBKK"

# ----- Row 43: KickToken -> CEJI -----
$ws.Cells.Item(43, 2).Value = "CEJI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 43 4 "0.002504"
$ws.Cells.Item(43, 5).Value = "42CEJICEJI"

# ----- Simple Price-only refreshes (continued) -----
Set-TextValue 44 4 "0.007854"
Set-TextValue 45 4 "0.00005334"

# ----- Row 48: BOLO - refreshed Price, suffix dropped from Volume -----
Set-TextValue 48 4 "0.002638"
$ws.Cells.Item(48, 5).Value = "47BOLOBOLO"

# Clean up the scratch cell/row so it doesn't leave a trace in the sheet.
$scratch.ClearContents()
$scratch.EntireRow.Delete()
